# "Generate Report for Handoff"
# Refresh the localization-status report: the two source files tracked in
# this run were handed off again (new GUID-prefixed filenames, new xlf
# build hashes, new timestamps), the zh-cn/de-de sheets lost their
# "Latest Target File" / "Latest Handback File" columns' data for this
# run (no handback happened yet), and the overview status flips from
# "Handed back: in sync with en-US" to "Ready for handoff".
#
# NOTE: this runtime's function calls don't bind named (-param value)
# arguments, so all helper functions below take positional args only.

$wb = $excel.ActiveWorkbook

# ---- helpers ---------------------------------------------------------

function Set-HyperlinkOnCell {
    param($ws, [string]$addr, [string]$target, [string]$display)
    $found = $null
    foreach ($h in $ws.Hyperlinks) {
        if ($h.Range.Address() -eq $addr) {
            $found = $h
        }
    }
    if ($found -ne $null) {
        $found.Address = $target
        $found.TextToDisplay = $display
    }
}

function Remove-HyperlinkOnCell {
    param($ws, [string]$addr)
    $found = $null
    foreach ($h in $ws.Hyperlinks) {
        if ($h.Range.Address() -eq $addr) {
            $found = $h
        }
    }
    if ($found -ne $null) {
        $found.Delete()
    }
}

# ---- new identifiers ---------------------------------------------------

$oldFile1 = "29e79f51-6ede-4853-a79d-4cea48aefdf7"
$oldFile2 = "69bc8315-b512-49ab-a3b9-5d471a9f1a0e"
$newFile1 = "1f641bcb-6eb4-4a1c-8351-6d8c6ffc8848"
$newFile2 = "ffff1f6804c2-9b89-4dcb-86c9-d98c33ec6e5e"

$newHash = "08af7258856ed5122eb31db14c922250c28858c4"

$newStatus = "Ready for handoff"
$newOverviewDate = "2016-47-18 08:47:19"

$newHandoffDateZh = "2016-03-18 08:47:16"
$newHandoffDateDe = "2016-03-18 08:47:19"
$clearedDate = "0001-01-01 00:00:00"

# =========================================================================
# Sheet "Overview"
# =========================================================================

$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "$newFile1.md"
$wsOverview.Range("B2").Value = $newStatus
$wsOverview.Range("C2").Value = $newStatus
$wsOverview.Range("D2").Value = $newOverviewDate

$wsOverview.Range("A3").Value = "$newFile2.md"
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus
$wsOverview.Range("D3").Value = $newOverviewDate

Set-HyperlinkOnCell $wsOverview '$A$2' "https://github.com/OpenLocalizationTest/oltest/blob/8e98f2051904923de05efd902c4af3697e35d9b1/e2e/$newFile1.md" "$newFile1.md"
Set-HyperlinkOnCell $wsOverview '$A$3' "https://github.com/OpenLocalizationTest/oltest/blob/8e98f2051904923de05efd902c4af3697e35d9b1/e2e/$newFile2.md" "$newFile2.md"

# =========================================================================
# Sheet "zh-cn"
# =========================================================================

$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A2").Value = "$newFile1.md"
$wsZh.Range("C2").Value = $newStatus
$wsZh.Range("D2").Value = "$newFile1.$newHash.zh-cn.xlf"
$wsZh.Range("E2").Value = $newHandoffDateZh
$wsZh.Range("H2").Value = $clearedDate

$wsZh.Range("A3").Value = "$newFile2.md"
$wsZh.Range("C3").Value = $newStatus
$wsZh.Range("D3").Value = "$newFile1.$newHash.zh-cn.xlf"
$wsZh.Range("E3").Value = $newHandoffDateZh
$wsZh.Range("H3").Value = $clearedDate

Set-HyperlinkOnCell $wsZh '$A$2' "https://github.com/OpenLocalizationTest/oltest/blob/8e98f2051904923de05efd902c4af3697e35d9b1/e2e/$newFile1.md" "$newFile1.md"
Set-HyperlinkOnCell $wsZh '$B$2' "https://github.com/OpenLocalizationTest/oltest/blob/8e98f2051904923de05efd902c4af3697e35d9b1/e2e/$newFile1.md" ".md"
Set-HyperlinkOnCell $wsZh '$D$2' "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a066214b83afa2279d9a41d28bdd90e0dc0a912e/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$newFile1.$newHash.zh-cn.xlf" "$newFile1.$newHash.zh-cn.xlf"

Set-HyperlinkOnCell $wsZh '$A$3' "https://github.com/OpenLocalizationTest/oltest/blob/8e98f2051904923de05efd902c4af3697e35d9b1/e2e/$newFile2.md" "$newFile2.md"
Set-HyperlinkOnCell $wsZh '$B$3' "https://github.com/OpenLocalizationTest/oltest/blob/8e98f2051904923de05efd902c4af3697e35d9b1/e2e/$newFile2.md" ".md"
Set-HyperlinkOnCell $wsZh '$D$3' "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a066214b83afa2279d9a41d28bdd90e0dc0a912e/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$newFile1.$newHash.zh-cn.xlf" "$newFile1.$newHash.zh-cn.xlf"

# "Latest Target File" (F) / "Latest Handback File" (G) have no data yet
# for this handoff run - drop their hyperlinks and wipe the cells
# outright (not just their contents) so the cell nodes disappear,
# matching the header-only F/G columns.
Remove-HyperlinkOnCell $wsZh '$F$2'
Remove-HyperlinkOnCell $wsZh '$G$2'
Remove-HyperlinkOnCell $wsZh '$F$3'
Remove-HyperlinkOnCell $wsZh '$G$3'
$wsZh.Range("F2:G3").Clear()

# =========================================================================
# Sheet "de-de"
# =========================================================================

$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A2").Value = "$newFile1.md"
$wsDe.Range("C2").Value = $newStatus
$wsDe.Range("D2").Value = "$newFile1.$newHash.de-de.xlf"
$wsDe.Range("E2").Value = $newHandoffDateDe
$wsDe.Range("H2").Value = $clearedDate

$wsDe.Range("A3").Value = "$newFile2.md"
$wsDe.Range("C3").Value = $newStatus
$wsDe.Range("D3").Value = "$newFile1.$newHash.de-de.xlf"
$wsDe.Range("E3").Value = $newHandoffDateDe
$wsDe.Range("H3").Value = $clearedDate

Set-HyperlinkOnCell $wsDe '$A$2' "https://github.com/OpenLocalizationTest/oltest/blob/8e98f2051904923de05efd902c4af3697e35d9b1/e2e/$newFile1.md" "$newFile1.md"
Set-HyperlinkOnCell $wsDe '$B$2' "https://github.com/OpenLocalizationTest/oltest/blob/8e98f2051904923de05efd902c4af3697e35d9b1/e2e/$newFile1.md" ".md"
Set-HyperlinkOnCell $wsDe '$D$2' "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8f2b64f41508255d80e8836c8296aa6cfd3860ac/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$newFile1.$newHash.de-de.xlf" "$newFile1.$newHash.de-de.xlf"

Set-HyperlinkOnCell $wsDe '$A$3' "https://github.com/OpenLocalizationTest/oltest/blob/8e98f2051904923de05efd902c4af3697e35d9b1/e2e/$newFile2.md" "$newFile2.md"
Set-HyperlinkOnCell $wsDe '$B$3' "https://github.com/OpenLocalizationTest/oltest/blob/8e98f2051904923de05efd902c4af3697e35d9b1/e2e/$newFile2.md" ".md"
Set-HyperlinkOnCell $wsDe '$D$3' "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8f2b64f41508255d80e8836c8296aa6cfd3860ac/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$newFile1.$newHash.de-de.xlf" "$newFile1.$newHash.de-de.xlf"

Remove-HyperlinkOnCell $wsDe '$F$2'
Remove-HyperlinkOnCell $wsDe '$G$2'
Remove-HyperlinkOnCell $wsDe '$F$3'
Remove-HyperlinkOnCell $wsDe '$G$3'
$wsDe.Range("F2:G3").Clear()
